# Add season-record columns (Wins / Losses / Ties) to the MIL_2023 sheet.
# Mirrors the commit "Created functions to get season record": the old
# scraper only pulled team statistics, not the won-loss-tied record, so
# three new columns are appended after the existing data (through AC).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting already used by the other header cells (bold,
# centered, bordered) by copying the format from the adjacent header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-60): every player row gets the same season record.
$ws.Range("AD2:AD60").Value = 92
$ws.Range("AE2:AE60").Value = 70
$ws.Range("AF2:AF60").Value = 0
